# Complete test plan for chequing account and update test values with
# accurate details (tests/A02_pixell_test_plan_chequing_account.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (C3, merged C3:D3)
$ws.Range("C3").Value = 'Sahil Choudhary'

# Row 7 - __init__ / balance greater than overdraft limit
$ws.Range("E7").Value = 'None'
$ws.Range("F7").Value = 'account_number = 350, client_number = 350, balance = 350'
$ws.Range("G7").Value = 'Attributes set to input values'

# Row 8 - __init__ / overdraft limit has invalid type.
$ws.Range("E8").Value = 'None'
$ws.Range("F8").Value = 'account_number = 350, client_number = 350, balance = 350, overdraft_limit = "-ten", overdraft_rate = 0.08'
$ws.Range("G8").Value = 'overdraft_limit set to -100'

# Row 9 - __init__ / overdraft rate has invalid type.
$ws.Range("E9").Value = 'None'
$ws.Range("F9").Value = 'account_number = 350, client_number = 350, balance = 350, overdraft_limit = -10, overdraft_rate = "eight percent"'
$ws.Range("G9").Value = 'overdraft_rate set to 0.05'

# Row 10 - __init__ / date created has invalid type
$ws.Range("E10").Value = 'None'
$ws.Range("F10").Value = 'account_number = 350, client_number = 350, balance = 350, date_created = "25 March, 2024", overdraft_limit = -10, overdraft_rate = 0.08'
$ws.Range("G10").Value = 'date_created set to current date'

# Row 11 - get_service_charges / balance greater than overdraft limit
$ws.Range("E11").Value = 'None'
$ws.Range("F11").Value = 'account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), overdraft_limit = -10, overdraft_rate = 0.08'
$ws.Range("G11").Value = 'service_charge = 0.50'

# Row 12 - get_service_charges / balance less than overdraft limit
$ws.Range("E12").Value = 'None'
$ws.Range("F12").Value = 'account_number = 350, client_number = 350, balance = -100, date_created = (2024, 3, 25), overdraft_limit = -10, overdraft_rate = 0.08'
$ws.Range("G12").Value = 'service_charge = 7.7'

# Row 13 - get_service_charges / balance equal to overdraft limit
$ws.Range("E13").Value = 'None'
$ws.Range("F13").Value = 'account_number = 350, client_number = 350, balance = -10, date_created = (2024, 3, 25), overdraft_limit = -10, overdraft_rate = 0.08'
$ws.Range("G13").Value = 'service_charge = 0.50'

# Row 14 - __str__ / appropriate value returned based on attribute values.
$ws.Range("E14").Value = 'None'
$ws.Range("F14").Value = 'account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), overdraft_limit = -10, overdraft_rate = 0.08'
$ws.Range("G14").Value = '"Account Number: 350 Balance: $350.00\nOverdraft Limit: $-10.00 Overdraft Rate: 8% Account Type: Chequing"'

# Mirror the final selection left by the editor (cell G14)
$ws.Range("G14").Select()
